$wb = $excel.ActiveWorkbook

# Mapping of row number -> new value for column F ("想去人数")
$updates = @{
    2  = 732
    3  = 599
    4  = 557
    7  = 85
    10 = 7
    11 = 4828
    12 = 4560
    13 = 6
    15 = 4
    16 = 35
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
